# Daily "Peru Liga 1" odds-feed refresh (30-05-2024 12:21).
#
# The upstream id sequence for match 7384628 was re-numbered: it now sorts
# into the slot currently held by id 7384625 (row 183), pushing the four
# rows that used to occupy 183-186 (ids 7384625, 7384626, 7384627, 7384630)
# down by one row, to 184-187. The running index in column A (181-185)
# stays put - it is the column B:AD match data (teams, odds, results, ...)
# that rotates up/down by one row.
#
# Net effect: row 187's B:AD content moves to row 183, and rows 183-186's
# B:AD content each shifts down into 184-187.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 183
$lastRow  = 187
$scratchRow = 1048500   # far outside the used range - safe holding cell

# 1) Stash the row that is moving to the top (old row 187) out of the way.
$ws.Range("B$lastRow`:AD$lastRow").Copy($ws.Range("B$scratchRow`:AD$scratchRow"))

# 2) Shift rows (lastRow-1)..firstRow down by one row each, bottom-up so we
#    never overwrite a source before it has been read.
for ($r = $lastRow - 1; $r -ge $firstRow; $r--) {
    $dst = $r + 1
    $ws.Range("B$r`:AD$r").Copy($ws.Range("B$dst`:AD$dst"))
}

# 3) Drop the stashed row into the now-vacated top slot and tidy up.
$ws.Range("B$scratchRow`:AD$scratchRow").Copy($ws.Range("B$firstRow`:AD$firstRow"))
$ws.Range("B$scratchRow`:AD$scratchRow").Clear()
